# Weekly refresh of the "Hortaliza, Agrícola del Norte S.A. de Arica - Berenjena"
# sheet: a new week's record is inserted at row 10, pushing the existing
# rows 10-19 down to rows 11-20 (dimension grows from A1:R19 to A1:R20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..19 down to 11..20, leaving row 10 free for the new record.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with this week's data.
$ws.Range("A10").Value2 = 1
$ws.Range("B10").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value2 = "Arica y Parinacota"
$ws.Range("D10").Value2 = 44827
$ws.Range("E10").Value2 = 15
$ws.Range("F10").Value2 = 100112001
$ws.Range("G10").Value2 = "Berenjena"
$ws.Range("H10").Value2 = "Sin especificar"
$ws.Range("I10").Value2 = "Primera"
$ws.Range("J10").Value2 = 120
$ws.Range("K10").Value2 = 6000
$ws.Range("L10").Value2 = 7000
$ws.Range("M10").Value2 = 6500
$ws.Range("N10").Value2 = "$/caja 60 unidades"
$ws.Range("O10").Value2 = "Región de Arica y Parinacota"
$ws.Range("P10").Value2 = 108
$ws.Range("Q10").Value2 = 60
$ws.Range("R10").Value2 = "Hortaliza"
